$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "BYD SEAL"
$ws.Range("C1").Value = "XPENG G9"
$ws.Range("D1").Value = "VINFAST VF8"
$ws.Range("E1").Value = "Honda ZR-V"
$ws.Range("F1").Value = "BYD SEAL-U"
$ws.Range("G1").Value = "Volkswagen ID.7"
$ws.Range("H1").Value = "BMW 5 series"
$ws.Range("I1").Value = "smart #3"
$ws.Range("J1").Value = "BYD TANG"
$ws.Range("K1").Value = "Hyundai KONA"
$ws.Range("L1").Value = "Kia EV9"
$ws.Range("M1").Value = "NIO ET5"
$ws.Range("N1").Value = "NIO EL7"
$ws.Range("O1").Value = "Lexus RZ"
